$wb = $excel.ActiveWorkbook

# --- PIR sheet: add rows 190-203 ---
$ws = $wb.Worksheets.Item("PIR")
$rng = $ws.Range("A190:F203")
$rng.NumberFormat = "@"

$ws.Cells.Item(190, 1).Value = "2026-01-28"
$ws.Cells.Item(190, 2).Value = "12:17:40"
$ws.Cells.Item(190, 3).Value = "12:00"
$ws.Cells.Item(190, 4).Value = "Bathroom"
$ws.Cells.Item(190, 5).Value = "No Motion"
$ws.Cells.Item(190, 6).Value = "Inactive"

$ws.Cells.Item(191, 1).Value = "2026-01-28"
$ws.Cells.Item(191, 2).Value = "12:17:40"
$ws.Cells.Item(191, 3).Value = "12:00"
$ws.Cells.Item(191, 4).Value = "Bathroom"
$ws.Cells.Item(191, 5).Value = "No Motion"
$ws.Cells.Item(191, 6).Value = "Inactive"

$ws.Cells.Item(192, 1).Value = "2026-01-28"
$ws.Cells.Item(192, 2).Value = "12:17:42"
$ws.Cells.Item(192, 3).Value = "12:00"
$ws.Cells.Item(192, 4).Value = "Bathroom"
$ws.Cells.Item(192, 5).Value = "No Motion"
$ws.Cells.Item(192, 6).Value = "Inactive"

$ws.Cells.Item(193, 1).Value = "2026-01-28"
$ws.Cells.Item(193, 2).Value = "12:17:46"
$ws.Cells.Item(193, 3).Value = "12:00"
$ws.Cells.Item(193, 4).Value = "Bathroom"
$ws.Cells.Item(193, 5).Value = "No Motion"
$ws.Cells.Item(193, 6).Value = "Inactive"

$ws.Cells.Item(194, 1).Value = "2026-01-28"
$ws.Cells.Item(194, 2).Value = "12:17:52"
$ws.Cells.Item(194, 3).Value = "12:00"
$ws.Cells.Item(194, 4).Value = "Bathroom"
$ws.Cells.Item(194, 5).Value = "No Motion"
$ws.Cells.Item(194, 6).Value = "Inactive"

$ws.Cells.Item(195, 1).Value = "2026-01-28"
$ws.Cells.Item(195, 2).Value = "12:17:56"
$ws.Cells.Item(195, 3).Value = "12:00"
$ws.Cells.Item(195, 4).Value = "Bathroom"
$ws.Cells.Item(195, 5).Value = "No Motion"
$ws.Cells.Item(195, 6).Value = "Inactive"

$ws.Cells.Item(196, 1).Value = "2026-01-28"
$ws.Cells.Item(196, 2).Value = "12:18:02"
$ws.Cells.Item(196, 3).Value = "12:00"
$ws.Cells.Item(196, 4).Value = "Bathroom"
$ws.Cells.Item(196, 5).Value = "No Motion"
$ws.Cells.Item(196, 6).Value = "Inactive"

$ws.Cells.Item(197, 1).Value = "2026-01-28"
$ws.Cells.Item(197, 2).Value = "12:18:07"
$ws.Cells.Item(197, 3).Value = "12:00"
$ws.Cells.Item(197, 4).Value = "Bathroom"
$ws.Cells.Item(197, 5).Value = "No Motion"
$ws.Cells.Item(197, 6).Value = "Inactive"

$ws.Cells.Item(198, 1).Value = "2026-01-28"
$ws.Cells.Item(198, 2).Value = "12:18:13"
$ws.Cells.Item(198, 3).Value = "12:00"
$ws.Cells.Item(198, 4).Value = "Bathroom"
$ws.Cells.Item(198, 5).Value = "No Motion"
$ws.Cells.Item(198, 6).Value = "Inactive"

$ws.Cells.Item(199, 1).Value = "2026-01-28"
$ws.Cells.Item(199, 2).Value = "12:18:17"
$ws.Cells.Item(199, 3).Value = "12:00"
$ws.Cells.Item(199, 4).Value = "Bathroom"
$ws.Cells.Item(199, 5).Value = "No Motion"
$ws.Cells.Item(199, 6).Value = "Inactive"

$ws.Cells.Item(200, 1).Value = "2026-01-28"
$ws.Cells.Item(200, 2).Value = "12:18:22"
$ws.Cells.Item(200, 3).Value = "12:00"
$ws.Cells.Item(200, 4).Value = "Bathroom"
$ws.Cells.Item(200, 5).Value = "No Motion"
$ws.Cells.Item(200, 6).Value = "Inactive"

$ws.Cells.Item(201, 1).Value = "2026-01-28"
$ws.Cells.Item(201, 2).Value = "12:18:27"
$ws.Cells.Item(201, 3).Value = "12:00"
$ws.Cells.Item(201, 4).Value = "Bathroom"
$ws.Cells.Item(201, 5).Value = "No Motion"
$ws.Cells.Item(201, 6).Value = "Inactive"

$ws.Cells.Item(202, 1).Value = "2026-01-28"
$ws.Cells.Item(202, 2).Value = "12:18:33"
$ws.Cells.Item(202, 3).Value = "12:00"
$ws.Cells.Item(202, 4).Value = "Bathroom"
$ws.Cells.Item(202, 5).Value = "No Motion"
$ws.Cells.Item(202, 6).Value = "Inactive"

$ws.Cells.Item(203, 1).Value = "2026-01-28"
$ws.Cells.Item(203, 2).Value = "12:18:37"
$ws.Cells.Item(203, 3).Value = "12:00"
$ws.Cells.Item(203, 4).Value = "Bathroom"
$ws.Cells.Item(203, 5).Value = "No Motion"
$ws.Cells.Item(203, 6).Value = "Inactive"

# --- Humidity sheet: add rows 179-190 ---
$ws = $wb.Worksheets.Item("Humidity")
$rng = $ws.Range("A179:F190")
$rng.NumberFormat = "@"

$ws.Cells.Item(179, 1).Value = "2026-01-28"
$ws.Cells.Item(179, 2).Value = "12:17:39"
$ws.Cells.Item(179, 3).Value = "12:00"
$ws.Cells.Item(179, 4).Value = "Bathroom"
$ws.Cells.Item(179, 5).Value = "87.8%"
$ws.Cells.Item(179, 6).Value = "Active"

$ws.Cells.Item(180, 1).Value = "2026-01-28"
$ws.Cells.Item(180, 2).Value = "12:17:41"
$ws.Cells.Item(180, 3).Value = "12:00"
$ws.Cells.Item(180, 4).Value = "Bathroom"
$ws.Cells.Item(180, 5).Value = "86.8%"
$ws.Cells.Item(180, 6).Value = "Active"

$ws.Cells.Item(181, 1).Value = "2026-01-28"
$ws.Cells.Item(181, 2).Value = "12:17:43"
$ws.Cells.Item(181, 3).Value = "12:00"
$ws.Cells.Item(181, 4).Value = "Bathroom"
$ws.Cells.Item(181, 5).Value = "86.2%"
$ws.Cells.Item(181, 6).Value = "Active"

$ws.Cells.Item(182, 1).Value = "2026-01-28"
$ws.Cells.Item(182, 2).Value = "12:17:47"
$ws.Cells.Item(182, 3).Value = "12:00"
$ws.Cells.Item(182, 4).Value = "Bathroom"
$ws.Cells.Item(182, 5).Value = "86.7%"
$ws.Cells.Item(182, 6).Value = "Active"

$ws.Cells.Item(183, 1).Value = "2026-01-28"
$ws.Cells.Item(183, 2).Value = "12:17:51"
$ws.Cells.Item(183, 3).Value = "12:00"
$ws.Cells.Item(183, 4).Value = "Bathroom"
$ws.Cells.Item(183, 5).Value = "87.6%"
$ws.Cells.Item(183, 6).Value = "Active"

$ws.Cells.Item(184, 1).Value = "2026-01-28"
$ws.Cells.Item(184, 2).Value = "12:18:03"
$ws.Cells.Item(184, 3).Value = "12:00"
$ws.Cells.Item(184, 4).Value = "Bathroom"
$ws.Cells.Item(184, 5).Value = "87.7%"
$ws.Cells.Item(184, 6).Value = "Active"

$ws.Cells.Item(185, 1).Value = "2026-01-28"
$ws.Cells.Item(185, 2).Value = "12:18:07"
$ws.Cells.Item(185, 3).Value = "12:00"
$ws.Cells.Item(185, 4).Value = "Bathroom"
$ws.Cells.Item(185, 5).Value = "86.7%"
$ws.Cells.Item(185, 6).Value = "Active"

$ws.Cells.Item(186, 1).Value = "2026-01-28"
$ws.Cells.Item(186, 2).Value = "12:18:12"
$ws.Cells.Item(186, 3).Value = "12:00"
$ws.Cells.Item(186, 4).Value = "Bathroom"
$ws.Cells.Item(186, 5).Value = "86.2%"
$ws.Cells.Item(186, 6).Value = "Active"

$ws.Cells.Item(187, 1).Value = "2026-01-28"
$ws.Cells.Item(187, 2).Value = "12:18:16"
$ws.Cells.Item(187, 3).Value = "12:00"
$ws.Cells.Item(187, 4).Value = "Bathroom"
$ws.Cells.Item(187, 5).Value = "87.7%"
$ws.Cells.Item(187, 6).Value = "Active"

$ws.Cells.Item(188, 1).Value = "2026-01-28"
$ws.Cells.Item(188, 2).Value = "12:18:20"
$ws.Cells.Item(188, 3).Value = "12:00"
$ws.Cells.Item(188, 4).Value = "Bathroom"
$ws.Cells.Item(188, 5).Value = "87.7%"
$ws.Cells.Item(188, 6).Value = "Active"

$ws.Cells.Item(189, 1).Value = "2026-01-28"
$ws.Cells.Item(189, 2).Value = "12:18:24"
$ws.Cells.Item(189, 3).Value = "12:00"
$ws.Cells.Item(189, 4).Value = "Bathroom"
$ws.Cells.Item(189, 5).Value = "87.6%"
$ws.Cells.Item(189, 6).Value = "Active"

$ws.Cells.Item(190, 1).Value = "2026-01-28"
$ws.Cells.Item(190, 2).Value = "12:18:32"
$ws.Cells.Item(190, 3).Value = "12:00"
$ws.Cells.Item(190, 4).Value = "Bathroom"
$ws.Cells.Item(190, 5).Value = "87.7%"
$ws.Cells.Item(190, 6).Value = "Active"

# --- Temperature sheet: add rows 179-190 ---
$ws = $wb.Worksheets.Item("Temperature")
$rng = $ws.Range("A179:F190")
$rng.NumberFormat = "@"

$ws.Cells.Item(179, 1).Value = "2026-01-28"
$ws.Cells.Item(179, 2).Value = "12:17:39"
$ws.Cells.Item(179, 3).Value = "12:00"
$ws.Cells.Item(179, 4).Value = "Bathroom"
$ws.Cells.Item(179, 5).Value = "23.0C"
$ws.Cells.Item(179, 6).Value = "Active"

$ws.Cells.Item(180, 1).Value = "2026-01-28"
$ws.Cells.Item(180, 2).Value = "12:17:42"
$ws.Cells.Item(180, 3).Value = "12:00"
$ws.Cells.Item(180, 4).Value = "Bathroom"
$ws.Cells.Item(180, 5).Value = "23.0C"
$ws.Cells.Item(180, 6).Value = "Active"

$ws.Cells.Item(181, 1).Value = "2026-01-28"
$ws.Cells.Item(181, 2).Value = "12:17:44"
$ws.Cells.Item(181, 3).Value = "12:00"
$ws.Cells.Item(181, 4).Value = "Bathroom"
$ws.Cells.Item(181, 5).Value = "23.0C"
$ws.Cells.Item(181, 6).Value = "Active"

$ws.Cells.Item(182, 1).Value = "2026-01-28"
$ws.Cells.Item(182, 2).Value = "12:17:48"
$ws.Cells.Item(182, 3).Value = "12:00"
$ws.Cells.Item(182, 4).Value = "Bathroom"
$ws.Cells.Item(182, 5).Value = "22.9C"
$ws.Cells.Item(182, 6).Value = "Active"

$ws.Cells.Item(183, 1).Value = "2026-01-28"
$ws.Cells.Item(183, 2).Value = "12:17:52"
$ws.Cells.Item(183, 3).Value = "12:00"
$ws.Cells.Item(183, 4).Value = "Bathroom"
$ws.Cells.Item(183, 5).Value = "23.0C"
$ws.Cells.Item(183, 6).Value = "Active"

$ws.Cells.Item(184, 1).Value = "2026-01-28"
$ws.Cells.Item(184, 2).Value = "12:18:04"
$ws.Cells.Item(184, 3).Value = "12:00"
$ws.Cells.Item(184, 4).Value = "Bathroom"
$ws.Cells.Item(184, 5).Value = "23.0C"
$ws.Cells.Item(184, 6).Value = "Active"

$ws.Cells.Item(185, 1).Value = "2026-01-28"
$ws.Cells.Item(185, 2).Value = "12:18:08"
$ws.Cells.Item(185, 3).Value = "12:00"
$ws.Cells.Item(185, 4).Value = "Bathroom"
$ws.Cells.Item(185, 5).Value = "22.9C"
$ws.Cells.Item(185, 6).Value = "Active"

$ws.Cells.Item(186, 1).Value = "2026-01-28"
$ws.Cells.Item(186, 2).Value = "12:18:12"
$ws.Cells.Item(186, 3).Value = "12:00"
$ws.Cells.Item(186, 4).Value = "Bathroom"
$ws.Cells.Item(186, 5).Value = "23.0C"
$ws.Cells.Item(186, 6).Value = "Active"

$ws.Cells.Item(187, 1).Value = "2026-01-28"
$ws.Cells.Item(187, 2).Value = "12:18:16"
$ws.Cells.Item(187, 3).Value = "12:00"
$ws.Cells.Item(187, 4).Value = "Bathroom"
$ws.Cells.Item(187, 5).Value = "23.0C"
$ws.Cells.Item(187, 6).Value = "Active"

$ws.Cells.Item(188, 1).Value = "2026-01-28"
$ws.Cells.Item(188, 2).Value = "12:18:20"
$ws.Cells.Item(188, 3).Value = "12:00"
$ws.Cells.Item(188, 4).Value = "Bathroom"
$ws.Cells.Item(188, 5).Value = "23.0C"
$ws.Cells.Item(188, 6).Value = "Active"

$ws.Cells.Item(189, 1).Value = "2026-01-28"
$ws.Cells.Item(189, 2).Value = "12:18:24"
$ws.Cells.Item(189, 3).Value = "12:00"
$ws.Cells.Item(189, 4).Value = "Bathroom"
$ws.Cells.Item(189, 5).Value = "22.9C"
$ws.Cells.Item(189, 6).Value = "Active"

$ws.Cells.Item(190, 1).Value = "2026-01-28"
$ws.Cells.Item(190, 2).Value = "12:18:32"
$ws.Cells.Item(190, 3).Value = "12:00"
$ws.Cells.Item(190, 4).Value = "Bathroom"
$ws.Cells.Item(190, 5).Value = "22.9C"
$ws.Cells.Item(190, 6).Value = "Active"
